# Updates cryptocurrency price/volume data in the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.192.48"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "2.422.34"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'554.45"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").Value = "'137.12"
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +1.98%  "
$ws.Range("E9").Value = "  -1.11%  "
$ws.Range("E10").Value = "  -0.60%  "
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("E12").Value = "  -1.15%  "
$ws.Range("E13").Value = "  -0.52%  "
$ws.Range("D14").Value = "2.855.70"
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("D15").Value = "60.103.77"
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("E16").Value = "  -0.33%  "
$ws.Range("D17").Value = "2.425.38"
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("D18").Value = "'11.26"
$ws.Range("E18").Value = "  -0.85%  "
$ws.Range("E19").Value = "  +2.88%  "
$ws.Range("D20").Value = "'328.02"
$ws.Range("E20").Value = "  -1.27%  "
$ws.Range("E21").Value = "  +0.65%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("E24").Value = "  +4.04%  "
$ws.Range("E25").Value = "  +0.73%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  +5.33%  "
$ws.Range("E28").Value = "  -1.19%  "
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("D30").Value = "'170.38"
$ws.Range("E30").Value = "  +0.23%  "
$ws.Range("E31").Value = "  -2.03%  "
$ws.Range("E32").Value = "  +3.41%  "
$ws.Range("E33").Value = "  -3.71%  "
$ws.Range("D34").Value = "'18.56"
$ws.Range("E34").Value = "  -0.57%  "
$ws.Range("E35").Value = "  +3.15%  "
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "'4.23"
$ws.Range("E37").Value = "  +0.26%  "
$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("D39").Value = "'325.60"
$ws.Range("E39").Value = "  +3.88%  "
$ws.Range("E40").Value = "  -0.79%  "
$ws.Range("D41").Value = "'146.15"
$ws.Range("E41").Value = "  +5.01%  "
$ws.Range("E43").Value = "  +0.20%  "
$ws.Range("D44").Value = "'19.79"
$ws.Range("E44").Value = "  +1.54%  "
$ws.Range("E45").Value = "  -0.60%  "
$ws.Range("E47").Value = "  -1.16%  "
$ws.Range("D48").Value = "'11.04"
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("E49").Value = "  -1.14%  "
$ws.Range("E50").Value = "  -0.54%  "
$ws.Range("D51").Value = "'0.945"
$ws.Range("E51").Value = "  -0.70%  "
